$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6542786666666667
$ws.Range("H2").Value = 1.962836
$ws.Range("I2").Value = 0.3193006097963691
$ws.Range("J2").Value = 0.3193006097963691
$ws.Range("M2").Value = 55.848606
$ws.Range("N2").Value = 167.545818
$ws.Range("O2").Value = 0.2323375192077237
$ws.Range("P2").Value = 0.2323375192077236
$ws.Range("Q2").Value = 36.540551468872
$ws.Range("R2").Value = 328.864963219848
$ws.Range("S2").Value = 0.07418551156160179
$ws.Range("T2").Value = 0.07418551156160177
$ws.Range("G3").Value = 0.6542786666666667
$ws.Range("H3").Value = 1.962836
$ws.Range("I3").Value = 0.3193006097963691
$ws.Range("J3").Value = 0.3193006097963691
$ws.Range("O3").Value = 0.3515710112922583
$ws.Range("P3").Value = 0.3515710112922583
$ws.Range("Q3").Value = 55.29282862662622
$ws.Range("R3").Value = 497.635457639636
$ws.Range("S3").Value = 0.1122568382923442
$ws.Range("T3").Value = 0.1122568382923442
$ws.Range("G4").Value = 0.6542786666666667
$ws.Range("H4").Value = 1.962836
$ws.Range("I4").Value = 0.3193006097963691
$ws.Range("J4").Value = 0.3193006097963691
$ws.Range("M4").Value = 33.195992
$ws.Range("N4").Value = 99.58797600000001
$ws.Range("O4").Value = 0.1380996766314891
$ws.Range("P4").Value = 0.1380996766314891
$ws.Range("Q4").Value = 21.71942938443734
$ws.Range("R4").Value = 195.474864459936
$ws.Range("S4").Value = 0.04409531096111586
$ws.Range("T4").Value = 0.04409531096111585
$ws.Range("G5").Value = 0.6542786666666667
$ws.Range("H5").Value = 1.962836
$ws.Range("I5").Value = 0.3193006097963691
$ws.Range("J5").Value = 0.3193006097963691
$ws.Range("M5").Value = 66.82284533333335
$ws.Range("N5").Value = 200.468536
$ws.Range("O5").Value = 0.277991792868529
$ws.Range("P5").Value = 0.2779917928685289
$ws.Range("Q5").Value = 43.72076214756623
$ws.Range("R5").Value = 393.486859328096
$ws.Range("S5").Value = 0.08876294898130724
$ws.Range("T5").Value = 0.08876294898130721
$ws.Range("I6").Value = 0.4124821994964292
$ws.Range("J6").Value = 0.4124821994964292
$ws.Range("M6").Value = 55.848606
$ws.Range("N6").Value = 167.545818
$ws.Range("O6").Value = 0.2323375192077237
$ws.Range("P6").Value = 0.2323375192077236
$ws.Range("Q6").Value = 47.204191217502
$ws.Range("R6").Value = 424.837720957518
$ws.Range("S6").Value = 0.09583509094834572
$ws.Range("T6").Value = 0.09583509094834571
$ws.Range("I7").Value = 0.4124821994964292
$ws.Range("J7").Value = 0.4124821994964292
$ws.Range("O7").Value = 0.3515710112922583
$ws.Range("P7").Value = 0.3515710112922583
$ws.Range("S7").Value = 0.1450167840170146
$ws.Range("T7").Value = 0.1450167840170146
$ws.Range("I8").Value = 0.4124821994964292
$ws.Range("J8").Value = 0.4124821994964292
$ws.Range("M8").Value = 33.195992
$ws.Range("N8").Value = 99.58797600000001
$ws.Range("O8").Value = 0.1380996766314891
$ws.Range("P8").Value = 0.1380996766314891
$ws.Range("Q8").Value = 28.057816770264
$ws.Range("R8").Value = 252.520350932376
$ws.Range("S8").Value = 0.05696365836670225
$ws.Range("T8").Value = 0.05696365836670224
$ws.Range("I9").Value = 0.4124821994964292
$ws.Range("J9").Value = 0.4124821994964292
$ws.Range("M9").Value = 66.82284533333335
$ws.Range("N9").Value = 200.468536
$ws.Range("O9").Value = 0.277991792868529
$ws.Range("P9").Value = 0.2779917928685289
$ws.Range("Q9").Value = 56.47980486410401
$ws.Range("R9").Value = 508.3182437769361
$ws.Range("S9").Value = 0.1146666661643666
$ws.Range("T9").Value = 0.1146666661643666
$ws.Range("G10").Value = 0.5380133333333333
$ws.Range("H10").Value = 1.61404
$ws.Range("I10").Value = 0.2625608844731457
$ws.Range("J10").Value = 0.2625608844731457
$ws.Range("M10").Value = 55.848606
$ws.Range("N10").Value = 167.545818
$ws.Range("O10").Value = 0.2323375192077237
$ws.Range("P10").Value = 0.2323375192077236
$ws.Range("Q10").Value = 30.04729467608
$ws.Range("R10").Value = 270.42565208472
$ws.Range("S10").Value = 0.06100274453947641
$ws.Range("T10").Value = 0.06100274453947641
$ws.Range("G11").Value = 0.5380133333333333
$ws.Range("H11").Value = 1.61404
$ws.Range("I11").Value = 0.2625608844731457
$ws.Range("J11").Value = 0.2625608844731457
$ws.Range("O11").Value = 0.3515710112922583
$ws.Range("P11").Value = 0.3515710112922583
$ws.Range("Q11").Value = 45.46729177400444
$ws.Range("R11").Value = 409.20562596604
$ws.Range("S11").Value = 0.09230879568001364
$ws.Range("T11").Value = 0.09230879568001364
$ws.Range("G12").Value = 0.5380133333333333
$ws.Range("H12").Value = 1.61404
$ws.Range("I12").Value = 0.2625608844731457
$ws.Range("J12").Value = 0.2625608844731457
$ws.Range("M12").Value = 33.195992
$ws.Range("N12").Value = 99.58797600000001
$ws.Range("O12").Value = 0.1380996766314891
$ws.Range("P12").Value = 0.1380996766314891
$ws.Range("Q12").Value = 17.85988630922667
$ws.Range("R12").Value = 160.73897678304
$ws.Range("S12").Value = 0.0362595732418192
$ws.Range("T12").Value = 0.03625957324181919
$ws.Range("G13").Value = 0.5380133333333333
$ws.Range("H13").Value = 1.61404
$ws.Range("I13").Value = 0.2625608844731457
$ws.Range("J13").Value = 0.2625608844731457
$ws.Range("M13").Value = 66.82284533333335
$ws.Range("N13").Value = 200.468536
$ws.Range("O13").Value = 0.277991792868529
$ws.Range("P13").Value = 0.2779917928685289
$ws.Range("Q13").Value = 35.95158176060445
$ws.Range("R13").Value = 323.56423584544
$ws.Range("S13").Value = 0.0729897710118365
$ws.Range("T13").Value = 0.07298977101183649
$ws.Range("G14").Value = 0.01159033333333333
$ws.Range("H14").Value = 0.034771
$ws.Range("I14").Value = 0.005656306234056004
$ws.Range("J14").Value = 0.005656306234056004
$ws.Range("M14").Value = 55.848606
$ws.Range("N14").Value = 167.545818
$ws.Range("O14").Value = 0.2323375192077237
$ws.Range("P14").Value = 0.2323375192077236
$ws.Range("Q14").Value = 0.6473039597420001
$ws.Range("R14").Value = 5.825735637678
$ws.Range("S14").Value = 0.001314172158299754
$ws.Range("T14").Value = 0.001314172158299754
$ws.Range("G15").Value = 0.01159033333333333
$ws.Range("H15").Value = 0.034771
$ws.Range("I15").Value = 0.005656306234056004
$ws.Range("J15").Value = 0.005656306234056004
$ws.Range("O15").Value = 0.3515710112922583
$ws.Range("P15").Value = 0.3515710112922583
$ws.Range("Q15").Value = 0.9794944377301111
$ws.Range("R15").Value = 8.815449939571
$ws.Range("S15").Value = 0.001988593302885774
$ws.Range("T15").Value = 0.001988593302885774
$ws.Range("G16").Value = 0.01159033333333333
$ws.Range("H16").Value = 0.034771
$ws.Range("I16").Value = 0.005656306234056004
$ws.Range("J16").Value = 0.005656306234056004
$ws.Range("M16").Value = 33.195992
$ws.Range("N16").Value = 99.58797600000001
$ws.Range("O16").Value = 0.1380996766314891
$ws.Range("P16").Value = 0.1380996766314891
$ws.Range("Q16").Value = 0.3847526126106667
$ws.Range("R16").Value = 3.462773513496001
$ws.Range("S16").Value = 0.0007811340618518101
$ws.Range("T16").Value = 0.00078113406185181
$ws.Range("G17").Value = 0.01159033333333333
$ws.Range("H17").Value = 0.034771
$ws.Range("I17").Value = 0.005656306234056004
$ws.Range("J17").Value = 0.005656306234056004
$ws.Range("M17").Value = 66.82284533333335
$ws.Range("N17").Value = 200.468536
$ws.Range("O17").Value = 0.277991792868529
$ws.Range("P17").Value = 0.2779917928685289
$ws.Range("Q17").Value = 0.7744990516951114
$ws.Range("R17").Value = 6.970491465256002
$ws.Range("S17").Value = 0.001572406711018666
$ws.Range("T17").Value = 0.001572406711018666
